# Auto-generated row-level updates for columns C (Median Value) and D (Tier)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; C = 0.7251761691223575; D = "Below Median" }
    @{ Row = 3; C = 0.9199231262011531; D = "Below Median" }
    @{ Row = 4; C = 2.65829596412556; D = $null }
    @{ Row = 5; C = 2.026905829596413; D = "2nd Tier" }
    @{ Row = 6; C = 3.136023916292975; D = $null }
    @{ Row = 7; C = 2.429212043561819; D = $null }
    @{ Row = 8; C = 1.596412556053812; D = "3rd Tier" }
    @{ Row = 9; C = 4.340807174887892; D = $null }
    @{ Row = 10; C = 0.6278026905829597; D = "Below Median" }
    @{ Row = 11; C = 0.7623318385650224; D = "Below Median" }
    @{ Row = 12; C = 0.8456117873158232; D = "Below Median" }
    @{ Row = 13; C = 0.9820627802690582; D = "Below Median" }
    @{ Row = 14; C = 0.9820627802690582; D = "Below Median" }
    @{ Row = 15; C = 0.8998505231689088; D = "Below Median" }
    @{ Row = 16; C = 1.086995515695067; D = "4th Tier" }
    @{ Row = 17; C = 1.051569506726457; D = "4th Tier" }
    @{ Row = 18; C = 1.237668161434978; D = "4th Tier" }
    @{ Row = 19; C = 0.9108121574489287; D = "Below Median" }
    @{ Row = 20; C = 0.7593423019431988; D = "Below Median" }
    @{ Row = 21; C = 0.7892376681614349; D = "Below Median" }
    @{ Row = 22; C = 0.8891736066623959; D = "Below Median" }
    @{ Row = 23; C = 0.8379244074311339; D = "Below Median" }
    @{ Row = 24; C = 2.895067264573991; D = $null }
    @{ Row = 25; C = 1.796284433055733; D = $null }
    @{ Row = 26; C = 3.019431988041854; D = $null }
    @{ Row = 27; C = 1.461883408071749; D = "3rd Tier" }
    @{ Row = 28; C = 2.125560538116592; D = "2nd Tier" }
    @{ Row = 29; C = 0.9650224215246637; D = "Below Median" }
    @{ Row = 30; C = 2.089686098654708; D = "2nd Tier" }
    @{ Row = 31; C = 1.285500747384156; D = "3rd Tier" }
    @{ Row = 32; C = 2.234977578475336; D = "2nd Tier" }
    @{ Row = 33; C = 0.5769805680119581; D = "Below Median" }
    @{ Row = 34; C = 2.078155028827675; D = "2nd Tier" }
    @{ Row = 35; C = 0.726457399103139; D = "Below Median" }
    @{ Row = 36; C = 2.208840486867393; D = "2nd Tier" }
    @{ Row = 37; C = 0.5739910313901345; D = "Below Median" }
    @{ Row = 38; C = 1.409352978859705; D = "3rd Tier" }
    @{ Row = 39; C = 1.194106342088405; D = "4th Tier" }
    @{ Row = 40; C = 2.787443946188341; D = $null }
    @{ Row = 41; C = 1; D = "4th Tier" }
    @{ Row = 42; C = 1.1898355754858; D = "4th Tier" }
    @{ Row = 43; C = 1.809417040358744; D = $null }
    @{ Row = 44; C = 1.456176110884631; D = "3rd Tier" }
    @{ Row = 45; C = 1.069058295964125; D = "4th Tier" }
    @{ Row = 46; C = 0.2654708520179372; D = $null }
    @{ Row = 47; C = 0.7294469357249626; D = "Below Median" }
    @{ Row = 48; C = 1.257847533632287; D = "3rd Tier" }
    @{ Row = 49; C = 0.6744394618834081; D = "Below Median" }
    @{ Row = 50; C = 0.9125560538116592; D = "Below Median" }
    @{ Row = 51; C = 0.2750373692077728; D = $null }
    @{ Row = 52; C = 0.8968609865470852; D = "Below Median" }
    @{ Row = 53; C = 2.517189835575486; D = $null }
    @{ Row = 54; C = 0.2331838565022422; D = $null }
    @{ Row = 55; C = 1.522101217168482; D = "3rd Tier" }
    @{ Row = 56; C = 0.9147982062780269; D = "Below Median" }
    @{ Row = 57; C = 1.280717488789238; D = "3rd Tier" }
    @{ Row = 58; C = 0.5704035874439461; D = "Below Median" }
    @{ Row = 59; C = 0.2600896860986547; D = $null }
    @{ Row = 60; C = 0.5944907110826394; D = "Below Median" }
    @{ Row = 61; C = 1.959641255605381; D = "2nd Tier" }
    @{ Row = 62; C = 2.946188340807175; D = $null }
    @{ Row = 63; C = 0.8143497757847533; D = "Below Median" }
    @{ Row = 64; C = 0.7461883408071749; D = "Below Median" }
    @{ Row = 65; C = 0.6816143497757847; D = "Below Median" }
    @{ Row = 66; C = 1.174887892376682; D = "4th Tier" }
    @{ Row = 67; C = 1.818834080717489; D = "2nd Tier" }
    @{ Row = 68; C = 0.8789237668161435; D = "Below Median" }
    @{ Row = 69; C = 0.2798206278026906; D = $null }
    @{ Row = 70; C = 0.7713004484304933; D = "Below Median" }
    @{ Row = 71; C = 0.5847533632286995; D = "Below Median" }
    @{ Row = 72; C = 1.356053811659193; D = "3rd Tier" }
    @{ Row = 73; C = 0.75695067264574; D = "Below Median" }
    @{ Row = 74; C = 1.191543882126842; D = "4th Tier" }
    @{ Row = 75; C = 4.573991031390134; D = $null }
    @{ Row = 76; C = 0.8938714499252616; D = "Below Median" }
    @{ Row = 77; C = 0.8968609865470852; D = "Below Median" }
    @{ Row = 78; C = 1.704035874439462; D = "3rd Tier" }
    @{ Row = 79; C = 0.57847533632287; D = "Below Median" }
    @{ Row = 80; C = 1.174887892376682; D = "4th Tier" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.C
    if ($u.D -ne $null) {
        $ws.Cells.Item($u.Row, 4).Value = $u.D
    }
}
